# Update ifo QoQ forecast error table (first evaluation) with refreshed
# preprocessing results for the ifo GDP component analysis.
# Only the data rows for quarters Q0-Q8 (rows 2-10) change; Q9 (row 11)
# stays as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    # row, ME, MAE, MSE, RMSE, SE, N
    @(2,  0.02037664140334729, 0.4125673874451535, 0.4288236258248655, 0.6548462612131687, 0.6609149592609375, 52),
    @(3,  0.4198709658825853,  0.6887531942785867, 1.402911000708223,  1.184445440156795,  1.118548952961554,  51),
    @(4,  0.5914175045278451,  1.023480694380445,  4.028787767966156,  2.007184039386064,  1.937548709935223,  50),
    @(5,  0.4886808719290308,  1.0759500315233,    4.644551171805161,  2.155122078167536,  2.120737718949353,  49),
    @(6,  0.4109263009013963,  0.9592852530966729, 4.304132647662862,  2.074640365861723,  2.055056383234047,  48),
    @(7,  0.4014302792999628,  1.035766140078607,  5.219666336070993,  2.284658910225111,  2.28101933610867,   36),
    @(8,  0.3912910724344761,  1.034972397376921,  5.326670609579055,  2.307958103947958,  2.307753342904947,  35),
    @(9,  0.2626343886760503,  1.460263562571626,  9.616065425151794,  3.100978140063518,  3.179415362967537,  18),
    @(10, -0.5588777053075298, 1.182169283809976,  6.428227015842427,  2.535394844169725,  2.593736879289778,  11)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
